$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 ("Комбинированный ввод" test): expected-result string gains "123"
# right before the trailing emoji.
$ws.Range("E10").Value = "TestВвода!123😊"

# Row 20 ("Копирование" test): expected-result text rewritten to describe the
# string being copied into the clipboard.
$ws.Range("F20").Value = "Строка копируется в буфер обмена"

# Update the sheet view: drop the scrolled topLeftCell and move the active
# selection from G25 to F21.
[void]$ws.Range("F21").Select()
